$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Kolkata_England): Matches_B 7->8, Won_B 4->5, Points_B 8->10
$ws.Range("H2").Value = 8
$ws.Range("I2").Value = 5
$ws.Range("K2").Value = 10

# Row 3 (Punjab_Pakistan): Matches_B 3->4, Lost_B 2->3
$ws.Range("H3").Value = 4
$ws.Range("J3").Value = 3

# Row 4 (Sunrisers_SriLanka): Matches_B 7->9, Won_B 3->4, Lost_B 4->5, Points_B 6->8
$ws.Range("H4").Value = 9
$ws.Range("I4").Value = 4
$ws.Range("J4").Value = 5
$ws.Range("K4").Value = 8
